$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value2 = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$ws.Range("B8").Value2 = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row after "Contact" (row 10),
# pushing Description/Purpose/... etc down by one row.
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row below (now row 12, formerly row 11)
# onto the newly inserted row 11 so it matches the sheet's existing style
# (border/alignment) without fabricating a brand-new style index.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value2 = "Jurisdiction"
$ws.Range("B11").Value2 = ""
